$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column B so Username/Password shift right
# (A: ScenarioName, B: Browsername [new], C: Username, D: Password)
$ws.Range("B1").EntireColumn.Insert()

# Header + data for the new "Browsername" column
$ws.Range("B1").Value = "Browsername"
$ws.Range("B2").Value = "chrome"
$ws.Range("B3").Value = "mozilla"

# Match the direct formatting (fill/border) of column A for the new column's cells
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)

$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# New column width (~24.14 characters)
$ws.Columns.Item(2).ColumnWidth = 23.3

# Match the selection of the target workbook
$ws.Range("D1").Select()
